$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header style from H1 onto I1:J1, then set header labels
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Fill in I0 and IF values for each data row (rows 2-55)
$iValues = @(7,8,6,9,11,8,9,6,6,6,7,8,10,6,8,8,4,9,7,7,7,6,5,8,9,8,8,8,1,7,5,7,4,8,2,7,8,6,8,3,8,8,7,9,7,7,7,7,4,8,7,4,3,6)
$jValues = @(7,9,7,9,11,8,9,6,6,6,7,8,10,6,8,8,5,9,7,7,7,7,5,8,9,8,8,8,2,7,6,8,5,9,3,7,9,7,9,3,8,9,8,9,9,7,7,7,5,8,7,4,3,6)

for ($idx = 0; $idx -lt $iValues.Length; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$idx]
    $ws.Cells.Item($row, 10).Value = $jValues[$idx]
}

